$p = $ppt.ActivePresentation
$s = $p.Slides.Add(4, 12)
